{"js": "// 1) Delete the whole paragraph:\n//    \"DENNE TESTEN B\u00d8R OGS\u00c5 GI VARSEL OM DET ER MER INFO, DA DETTE KAN V\u00c6RE TEGN P\u00c5 NOE GALT.\"\n{\n  const results = context.document.body.search(\n    \"DENNE TESTEN B\u00d8R OGS\u00c5 GI VARSEL OM DET ER MER INFO, DA DETTE KAN V\u00c6RE TEGN P\u00c5 NOE GALT.\",\n    { matchCase: true }\n  );\n  await context.sync();\n  if (results.items.length > 0) {\n    const para = results.items[0].paragraphs.getFirst();\n    para.delete();\n    await context.sync();\n  }\n}\n\n// 2) Move the <w:lastRenderedPageBreak/> marker: it currently sits on the \"AND/OR\" run\n//    that immediately follows the \"... en liten feil for \u00f8yeblikket ...\" sentence; it should\n//    instead sit at the start of the \"ANTALLDUPLISERTEKLASSER\" run (two paragraphs later).\n{\n  // Locate the correct \"AND/OR\" occurrence (there are several in the doc) by checking the\n  // paragraph two steps back for the unique anchor text.\n  const andOrResults = context.document.body.search(\"AND/OR\", { matchCase: true });\n  await context.sync();\n\n  let target = null;\n  for (const r of andOrResults.items) {\n    const para = r.paragraphs.getFirst();\n    const prev2 = para.getPrevious().getPrevious();\n    prev2.load(\"text\");\n    await context.sync();\n    if (prev2.text.indexOf(\"en liten feil\") !== -1) {\n      target = para;\n      break;\n    }\n  }\n\n  if (target) {\n    // Rewrite the \"AND/OR\" paragraph without the page-break marker.\n    const rng = target.getRange();\n    const plainOoxml =\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body><w:p><w:r><w:t>AND/OR</w:t></w:r></w:p></w:body>' +\n      '</w:document>' +\n      '</pkg:xmlData></pkg:part></pkg:package>';\n    rng.insertOoxml(plainOoxml, Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // Now add the marker to the start of the \"ANTALLDUPLISERTEKLASSER\" run.\n  const dupResults = context.document.body.search(\"ANTALLDUPLISERTEKLASSER\", { matchCase: true });\n  await context.sync();\n  if (dupResults.items.length > 0) {\n    const dupPara = dupResults.items[0].paragraphs.getFirst();\n    const dupRange = dupPara.getRange();\n    const dupOoxml =\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body><w:p>' +\n      '<w:r><w:rPr><w:b/><w:bCs/><w:u w:val=\"single\"/></w:rPr><w:lastRenderedPageBreak/><w:t>ANTALLDUPLISERTEKLASSER</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> klasser er dupliserte </w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\">i uttrekket. </w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\">Disse vises \\u00E5 v\\u00E6re dupliserte fra andre klassifikasjonssystem, og godkjennes. </w:t></w:r>' +\n      '<w:r><w:rPr><w:b/><w:bCs/><w:u w:val=\"single\"/></w:rPr><w:t>MANUELL SJEKK N\\u00C5R DE</w:t></w:r>' +\n      '<w:r><w:rPr><w:b/><w:bCs/><w:u w:val=\"single\"/></w:rPr><w:t>T SISTNEVNTE IKKE STEMMER.</w:t></w:r>' +\n      '</w:p></w:body>' +\n      '</w:document>' +\n      '</pkg:xmlData></pkg:part></pkg:package>';\n    dupRange.insertOoxml(dupOoxml, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 3) Merge the three runs \" \" + \"(referanseSekund\u00e6rKlassifikasjon\" + \")\" into a single run\n//    with text \" (referanseSekund\u00e6rKlassifikasjon)\".\n{\n  const results = context.document.body.search(\" (referanseSekund\u00e6rKlassifikasjon)\", { matchCase: true });\n  await context.sync();\n  if (results.items.length > 0) {\n    const rng = results.items[0];\n    rng.insertText(\" (referanseSekund\u00e6rKlassifikasjon)\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Delete the whole paragraph \"DENNE TESTEN B\u00d8R OGS\u00c5 GI VARSEL OM DET ER MER INFO, DA DETTE KAN V\u00c6RE TEGN P\u00c5 NOE GALT.\"\nforeach ($p in @($d.Paragraphs)) {\n    $t = $p.Range.Text\n    if ($t -match 'DENNE TESTEN B\u00d8R OGS\u00c5 GI VARSEL OM DET ER MER INFO, DA DETTE KAN V\u00c6RE TEGN P\u00c5 NOE GALT\\.') {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2) Move the lastRenderedPageBreak marker: remove it from the \"AND/OR\" run that precedes\n#    \"ANTALLDUPLISERTEKLASSER\" and add it at the start of the \"ANTALLDUPLISERTEKLASSER\" run.\n#    lastRenderedPageBreak is a rendering artifact; emulate by locating runs via Find.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"ANTALLDUPLISERTEKLASSER\"\nif ($rng.Find.Execute()) {\n    # no-op placeholder; actual marker handled at OOXML level below\n}\n\n# 3) Merge three runs \" \" + \"(referanseSekund\u00e6rKlassifikasjon\" + \")\" into one run text\n#    \" (referanseSekund\u00e6rKlassifikasjon)\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \" (referanseSekund\u00e6rKlassifikasjon)\"\nif (-not $rng2.Find.Execute()) {\n    $rng3 = $d.Content\n    $rng3.Find.ClearFormatting()\n    $rng3.Find.Text = \"(referanseSekund\u00e6rKlassifikasjon\"\n    $rng3.Find.Execute() | Out-Null\n}\n"}
